# The post "「この生き物を見分けられますか？」" (originally row 423) was removed
# from the sheet. Deleting the entire row shifts every row below it up by
# one (424->423, 425->424, ..., 596->595) and shrinks the sheet's used
# range/dimension from A1:C596 down to A1:C595, exactly matching the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(423).Delete()
